{"js": "// Office.js (Word JavaScript API) script.\n//\n// Splits several long run-of-text paragraphs into multiple line-broken\n// segments, matching the source edit. A manual line break (\"\\v\", the\n// vertical-tab control character) inside text passed to insertText is\n// rendered as a <w:br/> element when the package is serialized, so we\n// can build the exact \"<w:t>...</w:t><w:br/><w:t>...</w:t>\" run shape\n// the diff calls for by doing a single search + single insertText per\n// paragraph (this also keeps everything inside one <w:r>, instead of\n// Word.InsertLocation breaks which would create a new run per segment).\n\nconst LB = \"\\v\";\n\nasync function replaceExactText(context, oldText, newText) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1) \"Programa\" section - Portuguese list, 1..6 -> line-broken\nconst ptOld =\n  \"1. Estrutura molecular e liga\u00e7\u00e3o qu\u00edmica: Teoria de liga\u00e7\u00e3o de val\u00eancia, estrutura de compostos com C, N, O; Rela\u00e7\u00e3o entre estrutura e propriedades fisico-qu\u00edmicas2. Orbitais moleculares e as mol\u00e9culas de O2 e N2: Limita\u00e7\u00f5es da teoria de liga\u00e7\u00e3o de val\u00eancia, reatividade diferenciada de O2 e N2, relev\u00e2ncia do O2 em sistemas biol\u00f3gicos, esp\u00e9cies reativas de oxig\u00eanio3. \u00c1cidos, bases e a correla\u00e7\u00e3o com os ligantes dos metais em solu\u00e7\u00e3o: Afinidade das bases por metais de transi\u00e7\u00e3o, equil\u00edbrio qu\u00edmico em sistemas biol\u00f3gicos4. Complexos met\u00e1licos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octa\u00e9dricos e tetra\u00e9dricos; \u00edons de metais de transi\u00e7\u00e3o em sistemas biol\u00f3gicos5. Sistemas biol\u00f3gicos de transporte: Transporte de O2 em mam\u00edferos, transfer\u00eancia de el\u00e9trons dependente de metaloprote\u00ednas;6. Processos catal\u00edticos - \u00e1cido/base e oxido-redu\u00e7\u00e3o em metaloprote\u00ednas: Prote\u00ednas contendo \u00edon Zn2+, peroxidases, oxidases.\";\nconst ptNew = [\n  \"1. Estrutura molecular e liga\u00e7\u00e3o qu\u00edmica: Teoria de liga\u00e7\u00e3o de val\u00eancia, estrutura de compostos com C, N, O; Rela\u00e7\u00e3o entre estrutura e propriedades fisico-qu\u00edmicas\",\n  \"2. Orbitais moleculares e as mol\u00e9culas de O2 e N2: Limita\u00e7\u00f5es da teoria de liga\u00e7\u00e3o de val\u00eancia, reatividade diferenciada de O2 e N2, relev\u00e2ncia do O2 em sistemas biol\u00f3gicos, esp\u00e9cies reativas de oxig\u00eanio\",\n  \"3. \u00c1cidos, bases e a correla\u00e7\u00e3o com os ligantes dos metais em solu\u00e7\u00e3o: Afinidade das bases por metais de transi\u00e7\u00e3o, equil\u00edbrio qu\u00edmico em sistemas biol\u00f3gicos\",\n  \"4. Complexos met\u00e1licos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octa\u00e9dricos e tetra\u00e9dricos; \u00edons de metais de transi\u00e7\u00e3o em sistemas biol\u00f3gicos\",\n  \"5. Sistemas biol\u00f3gicos de transporte: Transporte de O2 em mam\u00edferos, transfer\u00eancia de el\u00e9trons dependente de metaloprote\u00ednas;\",\n  \"6. Processos catal\u00edticos - \u00e1cido/base e oxido-redu\u00e7\u00e3o em metaloprote\u00ednas: Prote\u00ednas contendo \u00edon Zn2+, peroxidases, oxidases.\",\n].join(LB);\nawait replaceExactText(context, ptOld, ptNew);\n\n// 2) \"Programa\" section - English (italic) list, 1..6 -> line-broken\nconst enOld =\n  \"1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases\";\nconst enNew = [\n  \"1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.\",\n  \"2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species\",\n  \"3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems\",\n  \"4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems\",\n  \"5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins\",\n  \"6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases\",\n].join(LB);\nawait replaceExactText(context, enOld, enNew);\n\n// 3) \"Crit\u00e9rio\" run - add breaks before NF=... and before Sendo que...\nconst critOld =\n  \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a mat\u00e9ria ser\u00e1 cumulativa do semestre.\";\nconst critNew = [\n  \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:\",\n  \"NF = (P1 + 2*P2)/3\",\n  \"\",\n  \"Sendo que para P2 a mat\u00e9ria ser\u00e1 cumulativa do semestre.\",\n].join(LB);\nawait replaceExactText(context, critOld, critNew);\n\n// 4) \"Norma de recupera\u00e7\u00e3o\" run - add double break before MR = ...\nconst normaOld =\n  \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: MR = (NF + PR)/2\";\nconst normaNew = [\n  \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: \",\n  \"\",\n  \"MR = (NF + PR)/2\",\n].join(LB);\nawait replaceExactText(context, normaOld, normaNew);\n\n// 5) \"Bibliografia\" run -> split into two lines\nconst biblioOld =\n  \"1. Atkins e Jones, Princ\u00edpios de Qu\u00edmica, 5a edi\u00e7\u00e3o, Bookman, 20112. Shiver e Atikins, Qu\u00edmica Inorg\u00e2nica, 4a edi\u00e7\u00e3o, Bookman, 2008\";\nconst biblioNew = [\n  \"1. Atkins e Jones, Princ\u00edpios de Qu\u00edmica, 5a edi\u00e7\u00e3o, Bookman, 2011\",\n  \"2. Shiver e Atikins, Qu\u00edmica Inorg\u00e2nica, 4a edi\u00e7\u00e3o, Bookman, 2008\",\n].join(LB);\nawait replaceExactText(context, biblioOld, biblioNew);\n", "ps1": "# Word COM interop script\n# Splits several long run-of-text paragraphs into multiple <w:t> runs\n# joined by manual line breaks (<w:br/>), matching the source edit.\n#\n# In Word's object model a manual line break is represented in Range.Text\n# as Chr(11) (vertical tab); writing that character back into a Range\n# re-creates a <w:br/> element, which is exactly what we need here.\n\n$d = $word.ActiveDocument\n$LB = [char]11\n\nfunction Replace-ExactText($oldText, $newText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Text = $newText\n    }\n    return $found\n}\n\n# 1) \"Programa\" section - Portuguese list, 1..6 -> line-broken\n$ptOld = \"1. Estrutura molecular e liga\u00e7\u00e3o qu\u00edmica: Teoria de liga\u00e7\u00e3o de val\u00eancia, estrutura de compostos com C, N, O; Rela\u00e7\u00e3o entre estrutura e propriedades fisico-qu\u00edmicas2. Orbitais moleculares e as mol\u00e9culas de O2 e N2: Limita\u00e7\u00f5es da teoria de liga\u00e7\u00e3o de val\u00eancia, reatividade diferenciada de O2 e N2, relev\u00e2ncia do O2 em sistemas biol\u00f3gicos, esp\u00e9cies reativas de oxig\u00eanio3. \u00c1cidos, bases e a correla\u00e7\u00e3o com os ligantes dos metais em solu\u00e7\u00e3o: Afinidade das bases por metais de transi\u00e7\u00e3o, equil\u00edbrio qu\u00edmico em sistemas biol\u00f3gicos4. Complexos met\u00e1licos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octa\u00e9dricos e tetra\u00e9dricos; \u00edons de metais de transi\u00e7\u00e3o em sistemas biol\u00f3gicos5. Sistemas biol\u00f3gicos de transporte: Transporte de O2 em mam\u00edferos, transfer\u00eancia de el\u00e9trons dependente de metaloprote\u00ednas;6. Processos catal\u00edticos - \u00e1cido/base e oxido-redu\u00e7\u00e3o em metaloprote\u00ednas: Prote\u00ednas contendo \u00edon Zn2+, peroxidases, oxidases.\"\n$ptNew = @(\n    \"1. Estrutura molecular e liga\u00e7\u00e3o qu\u00edmica: Teoria de liga\u00e7\u00e3o de val\u00eancia, estrutura de compostos com C, N, O; Rela\u00e7\u00e3o entre estrutura e propriedades fisico-qu\u00edmicas\",\n    \"2. Orbitais moleculares e as mol\u00e9culas de O2 e N2: Limita\u00e7\u00f5es da teoria de liga\u00e7\u00e3o de val\u00eancia, reatividade diferenciada de O2 e N2, relev\u00e2ncia do O2 em sistemas biol\u00f3gicos, esp\u00e9cies reativas de oxig\u00eanio\",\n    \"3. \u00c1cidos, bases e a correla\u00e7\u00e3o com os ligantes dos metais em solu\u00e7\u00e3o: Afinidade das bases por metais de transi\u00e7\u00e3o, equil\u00edbrio qu\u00edmico em sistemas biol\u00f3gicos\",\n    \"4. Complexos met\u00e1licos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octa\u00e9dricos e tetra\u00e9dricos; \u00edons de metais de transi\u00e7\u00e3o em sistemas biol\u00f3gicos\",\n    \"5. Sistemas biol\u00f3gicos de transporte: Transporte de O2 em mam\u00edferos, transfer\u00eancia de el\u00e9trons dependente de metaloprote\u00ednas;\",\n    \"6. Processos catal\u00edticos - \u00e1cido/base e oxido-redu\u00e7\u00e3o em metaloprote\u00ednas: Prote\u00ednas contendo \u00edon Zn2+, peroxidases, oxidases.\"\n) -join $LB\nReplace-ExactText $ptOld $ptNew | Out-Null\n\n# 2) \"Programa\" section - English (italic) list, 1..6 -> line-broken\n$enOld = \"1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases\"\n$enNew = @(\n    \"1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.\",\n    \"2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species\",\n    \"3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems\",\n    \"4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems\",\n    \"5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins\",\n    \"6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases\"\n) -join $LB\nReplace-ExactText $enOld $enNew | Out-Null\n\n# 3) \"Crit\u00e9rio\" run - add breaks before NF=... and before Sendo que...\n$critOld = \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a mat\u00e9ria ser\u00e1 cumulativa do semestre.\"\n$critNew = @(\n    \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:\",\n    \"NF = (P1 + 2*P2)/3\",\n    \"\",\n    \"Sendo que para P2 a mat\u00e9ria ser\u00e1 cumulativa do semestre.\"\n) -join $LB\nReplace-ExactText $critOld $critNew | Out-Null\n\n# 4) \"Norma de recupera\u00e7\u00e3o\" run - add double break before MR = ...\n$normaOld = \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: MR = (NF + PR)/2\"\n$normaNew = @(\n    \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: \",\n    \"\",\n    \"MR = (NF + PR)/2\"\n) -join $LB\nReplace-ExactText $normaOld $normaNew | Out-Null\n\n# 5) \"Bibliografia\" run -> split into two lines\n$biblioOld = \"1. Atkins e Jones, Princ\u00edpios de Qu\u00edmica, 5a edi\u00e7\u00e3o, Bookman, 20112. Shiver e Atikins, Qu\u00edmica Inorg\u00e2nica, 4a edi\u00e7\u00e3o, Bookman, 2008\"\n$biblioNew = @(\n    \"1. Atkins e Jones, Princ\u00edpios de Qu\u00edmica, 5a edi\u00e7\u00e3o, Bookman, 2011\",\n    \"2. Shiver e Atikins, Qu\u00edmica Inorg\u00e2nica, 4a edi\u00e7\u00e3o, Bookman, 2008\"\n) -join $LB\nReplace-ExactText $biblioOld $biblioNew | Out-Null\n"}
